# Append: 2025-09-04 18:31 JST
# - Refresh the "取得日時" (fetched-at) timestamp for every data row (2-15)
#   from 2025-09-04 18:25:04 to 2025-09-04 18:31:29.
# - Row 13 and row 14 swap their listing (タイトル/URL), matching the
#   site re-ordering the "注目" and "限定公開" PR slots.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-09-04 18:25:04"
$newTimestamp = "2025-09-04 18:31:29"

# Update the fetch timestamp in column A for rows 2 through 15.
for ($r = 2; $r -le 15; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cur = $cell.Value2
    if ($cur -eq $oldTimestamp) {
        $cell.Value2 = $newTimestamp
    }
}

# Rows 13 and 14 swap their title (B) and URL (F) contents.
$b13 = $ws.Range("B13").Value2
$f13 = $ws.Range("F13").Value2
$b14 = $ws.Range("B14").Value2
$f14 = $ws.Range("F14").Value2

$ws.Range("B13").Value2 = $b14
$ws.Range("F13").Value2 = $f14
$ws.Range("B14").Value2 = $b13
$ws.Range("F14").Value2 = $f13
